# Bug fix in Eduati data files (SW837_noCTRL_meas.xlsx)
#
# 1. Sheet1 ("measurements") had 86 data rows (rows 2:87) but only the
#    first 43 (rows 2:44) actually carry real B:N measurement values -
#    rows 45:87 are stray leftover index-only rows (column A only) from
#    a previous, longer dataset. Delete them so Sheet1 lines up with
#    Sheet2/Sheet3 (both A1:N44).
# 2. Re-point the active tab/selection: Sheet1 becomes the active sheet
#    (scrolled near the bottom of the now-44-row table, cell C56 - which
#    is below the data - ends up selected), and Sheet3 is no longer the
#    tab that was left selected.
#
# Note: the underlying x15ac:absPath (last-saved-folder breadcrumb) is
# Microsoft-managed metadata that Excel itself stamps on save from the
# real filesystem path in use - it is not part of the Workbook/Worksheet
# object model and has no COM/VBA setter, so it cannot be scripted here;
# it is left untouched by this automation.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Drop the 43 stray index-only rows (45:87) at the bottom of Sheet1 -
# this also shrinks the sheet's dimension from A1:N87 down to A1:N44.
$ws1.Rows("45:87").Delete() | Out-Null

# Sheet1 becomes the active/selected sheet and tab (previously Sheet3
# was last active via workbook-level activeTab); select C56 on it.
$ws1.Activate() | Out-Null
$ws1.Range("C56").Select() | Out-Null

# Sheet3's own in-sheet selection (A2:N44) is left as-is; only its
# "this was the active tab" flag goes away, which happens automatically
# now that Sheet1 is activated above.
